$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()
$ws.Range("H98").Value = 2203.3125
$ws.Range("I98").Value = 2058
$ws.Range("K98").Value = 2058
$ws.Range("M98").Value = -560
$ws.Range("H100").Value = 3428.8823
$ws.Range("I100").Value = 1656.2142
$ws.Range("J100").Value = 4669.75
$ws.Range("K100").Value = 1656.2142
$ws.Range("L100").Value = 4669.75
$ws.Range("M100").Value = -1115.2142
$ws.Range("N100").Value = -5751.75
$ws.Range("H122").Value = 2203.3125
$ws.Range("I122").Value = 2058
$ws.Range("K122").Value = 6174
$ws.Range("M122").Value = -3724
$ws.Range("H139").Value = 74865
$ws.Range("J139").Value = 74865
$ws.Range("L139").Value = 74865
$ws.Range("N139").Value = -85145
$ws.Range("H141").Value = 7760.3887
$ws.Range("I141").Value = 4326.636
$ws.Range("K141").Value = 12979.908
$ws.Range("M141").Value = -7799.908000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H44").Value = 26969.8
$ws.Range("J44").Value = 26969.8
$ws.Range("L44").Value = 26969.8
$ws.Range("N44").Value = -27945.8
$ws.Range("H55").Value = 29600
$ws.Range("J55").Value = 29600
$ws.Range("L55").Value = 29600
$ws.Range("N55").Value = -30230
$ws.Range("H63").Value = 4160.5454
$ws.Range("I63").Value = 2148.8333
$ws.Range("J63").Value = 6574.6
$ws.Range("K63").Value = 2148.8333
$ws.Range("L63").Value = 6574.6
$ws.Range("M63").Value = -1462.8333
$ws.Range("N63").Value = -7946.6
$ws.Range("H66").Value = 4160.5454
$ws.Range("I66").Value = 2148.8333
$ws.Range("J66").Value = 6574.6
$ws.Range("K66").Value = 10744.1665
$ws.Range("L66").Value = 32873
$ws.Range("M66").Value = -7312.166499999999
$ws.Range("N66").Value = -39737
$ws.Range("H80").Value = 39544.547
$ws.Range("J80").Value = 39499
$ws.Range("L80").Value = 39499
$ws.Range("N80").Value = -41495
$ws.Range("H83").Value = 39544.547
$ws.Range("J83").Value = 39499
$ws.Range("L83").Value = 118497
$ws.Range("N83").Value = -128481
$ws.Range("H88").Value = 1653.8125
$ws.Range("J88").Value = 1878.8182
$ws.Range("L88").Value = 1878.8182
$ws.Range("N88").Value = -2690.8182
$ws.Range("H91").Value = 1653.8125
$ws.Range("J91").Value = 1878.8182
$ws.Range("L91").Value = 1878.8182
$ws.Range("N91").Value = -4686.8182
$ws.Range("H110").Value = 1418.909
$ws.Range("I110").Value = 1053.7
$ws.Range("J110").Value = 1980.7693
$ws.Range("K110").Value = 1053.7
$ws.Range("L110").Value = 1980.7693
$ws.Range("M110").Value = 991.3
$ws.Range("N110").Value = -6070.7693
$ws.Range("H134").Value = 60000.668
$ws.Range("J134").Value = 60000.668
$ws.Range("L134").Value = 60000.668
$ws.Range("N134").Value = -70140.66800000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 32704.2
$ws.Range("I35").Value = 25000
$ws.Range("J35").Value = 34630.25
$ws.Range("K35").Value = 25000
$ws.Range("L35").Value = 34630.25
$ws.Range("M35").Value = -24690
$ws.Range("N35").Value = -35250.25
$ws.Range("H82").Value = 27776.334
$ws.Range("J82").Value = 43455.668
$ws.Range("L82").Value = 43455.668
$ws.Range("N82").Value = -44221.668
$ws.Range("H85").Value = 27776.334
$ws.Range("J85").Value = 43455.668
$ws.Range("L85").Value = 43455.668
$ws.Range("N85").Value = -46107.668
$ws.Range("H94").Value = 3375.9565
$ws.Range("I94").Value = 3554.6191
$ws.Range("K94").Value = 3554.6191
$ws.Range("M94").Value = -3103.6191
$ws.Range("H107").Value = 1579.6957
$ws.Range("I107").Value = 1256.3846
$ws.Range("K107").Value = 1256.3846
$ws.Range("M107").Value = 663.6153999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 60752.05
$ws.Range("I16").Value = 67717.64999999999
$ws.Range("J16").Value = 1544.5
$ws.Range("K16").Value = 67717.64999999999
$ws.Range("L16").Value = 1544.5
$ws.Range("M16").Value = -67430.64999999999
$ws.Range("N16").Value = -2118.5
$ws.Range("H19").Value = 2778407.2
$ws.Range("I19").Value = 3572194.2
$ws.Range("J19").Value = 152.75
$ws.Range("K19").Value = 3572194.2
$ws.Range("L19").Value = 152.75
$ws.Range("M19").Value = -3572024.2
$ws.Range("N19").Value = -492.75
$ws.Range("H24").Value = 2778407.2
$ws.Range("I24").Value = 3572194.2
$ws.Range("J24").Value = 152.75
$ws.Range("K24").Value = 3572194.2
$ws.Range("L24").Value = 152.75
$ws.Range("M24").Value = -3572024.2
$ws.Range("N24").Value = -492.75
$ws.Range("H86").Value = 38304.9
$ws.Range("I86").Value = 94528.664
$ws.Range("J86").Value = 13004.2
$ws.Range("K86").Value = 94528.664
$ws.Range("L86").Value = 13004.2
$ws.Range("M86").Value = -93405.664
$ws.Range("N86").Value = -15250.2
$ws.Range("H89").Value = 38304.9
$ws.Range("I89").Value = 94528.664
$ws.Range("J89").Value = 13004.2
$ws.Range("K89").Value = 472643.32
$ws.Range("L89").Value = 65021
$ws.Range("M89").Value = -467027.32
$ws.Range("N89").Value = -76253
$ws.Range("H99").Value = 23333.8
$ws.Range("I99").Value = 31706.857
$ws.Range("J99").Value = 3796.6667
$ws.Range("K99").Value = 31706.857
$ws.Range("L99").Value = 3796.6667
$ws.Range("M99").Value = -30208.857
$ws.Range("N99").Value = -6792.6667
$ws.Range("H107").Value = 599.43475
$ws.Range("I107").Value = 442.2857
$ws.Range("K107").Value = 442.2857
$ws.Range("M107").Value = 1477.7143
$ws.Range("H113").Value = 60752.05
$ws.Range("I113").Value = 67717.64999999999
$ws.Range("J113").Value = 1544.5
$ws.Range("K113").Value = 67717.64999999999
$ws.Range("L113").Value = 1544.5
$ws.Range("M113").Value = -65547.64999999999
$ws.Range("N113").Value = -5884.5
$ws.Range("H126").Value = 23333.8
$ws.Range("I126").Value = 31706.857
$ws.Range("J126").Value = 3796.6667
$ws.Range("K126").Value = 95120.571
$ws.Range("L126").Value = 11390.0001
$ws.Range("M126").Value = -92650.571
$ws.Range("N126").Value = -16330.0001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 5041250
$ws.Range("I11").Value = 5041250
$ws.Range("K11").Value = 5041250
$ws.Range("M11").Value = -5041111
$ws.Range("H22").Value = 500
$ws.Range("J22").Value = 500
$ws.Range("L22").Value = 500
$ws.Range("N22").Value = -1558
$ws.Range("H46").Value = 19000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 19000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 19000
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -19312
$ws.Range("H57").Value = 19771.8
$ws.Range("J57").Value = 22214.75
$ws.Range("L57").Value = 22214.75
$ws.Range("N57").Value = -23854.75
$ws.Range("H80").Value = 1313.5
$ws.Range("J80").Value = 1599.5
$ws.Range("L80").Value = 1599.5
$ws.Range("N80").Value = -3595.5
$ws.Range("H83").Value = 1313.5
$ws.Range("J83").Value = 1599.5
$ws.Range("L83").Value = 7997.5
$ws.Range("N83").Value = -17981.5
$ws.Range("H97").Value = 1256.4572
$ws.Range("I97").Value = 1197.1538
$ws.Range("K97").Value = 1197.1538
$ws.Range("M97").Value = -701.1538
$ws.Range("H113").Value = 1650
$ws.Range("I113").Value = 1650
$ws.Range("K113").Value = 1650
$ws.Range("M113").Value = 520

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1286486.8
$ws.Range("I132").Value = 1854550
$ws.Range("K132").Value = 5563650
$ws.Range("M132").Value = -5561120

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 6668751
$ws.Range("I132").Value = 7577790
$ws.Range("K132").Value = 22733370
$ws.Range("M132").Value = -22730840
$ws.Range("H136").Value = 7642220.5
$ws.Range("I136").Value = 4348889
$ws.Range("K136").Value = 13046667
$ws.Range("M136").Value = -13044117

Write-Host "All changes applied"